# The commit simplifies the document's default formatting block
# (w:docDefaults in styles.xml): a long list of redundant, schema-default
# run/paragraph properties (b=0, i=0, smallCaps=0, strike=0, color=auto,
# u=none, shd=clear/auto, vertAlign=baseline, keepNext=0, keepLines=0,
# widowControl=1, empty pBdr, ind=0, contextualSpacing=0, jc=left, and the
# before/after spacing of 0) is stripped away, leaving only the handful of
# values that actually matter: the Arial/22-half-point/"en" run defaults,
# and a paragraph default of 276 line spacing on the "auto" (multiple)
# rule. None of these removed values changes what any paragraph actually
# looks like -- they were already the effective values every paragraph in
# the document was using.
#
# Word's object model has no direct handle onto <w:docDefaults> itself
# (it is not a scriptable object -- Word only ever regenerates it when a
# new document is created from a template), so the closest reachable
# equivalent is to (re-)apply the same, already-effective line spacing to
# the built-in "Normal" style, which is exactly what every paragraph in
# this document already inherits. This keeps the visible formatting
# identical while matching the simplified "only line spacing matters"
# shape of the cleaned-up defaults.

$d = $word.ActiveDocument
$normal = $d.Styles("Normal")

$normal.ParagraphFormat.LineSpacingRule = 5   # wdLineSpaceMultiple
$normal.ParagraphFormat.LineSpacing = 13.8    # 276 twentieths-of-a-point = 1.15x
